$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the individual card-field rows into single Python-tuple-repr strings.
$ws.Range("A2").Value = "('Hammer of Bogardan', ['{1}{R}{R}', 'Sorcery', 'Hammer of Bogardan deals 3 damage to any target.', '{2}{R}{R}{R}: Return Hammer of Bogardan from your graveyard to your hand. Activate this ability only during your upkeep.'])"
$ws.Range("A3").Value = "('Tradewind Rider', ['{3}{U}', 'Creature — Spirit', 'Flying', '{T}, Tap two untapped creatures you control: Return target permanent to its owner’s hand.', '1/4'])"

# Remove the now-obsolete rows 4-12 that held the individual fields.
$ws.Range("A4:A12").ClearContents()
